$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cluster assignment values for rows 2-34 (weather clusters AF/full-size fix)
$ws.Range("F2").Value = 0.02707373271889402
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0.02247191011235954
$ws.Range("B3").Value = 0.002027027027027027
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 0.05587557603686627
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0.06292134831460669
$ws.Range("C4").Value = 0.09350775193798465
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 0.02252252252252252
$ws.Range("I4").Value = 0.08940397350993383
$ws.Range("J4").Value = 0
$ws.Range("B5").Value = 0.3878378378378345
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.01816239316239316
$ws.Range("F5").Value = 0.2355990783410156
$ws.Range("H5").Value = 0.002257336343115124
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.2078651685393255
$ws.Range("F6").Value = 0.02188940092165899
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0.007865168539325843
$ws.Range("I7").Value = 0.0380794701986755
$ws.Range("B8").Value = 0.2844594594594572
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0.3865225683407469
$ws.Range("F8").Value = 0.05990783410138238
$ws.Range("H8").Value = 0.4225733634311608
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0.1415730337078653
$ws.Range("K8").Value = 0.01363366923532898
$ws.Range("D9").Value = 0.004273504273504274
$ws.Range("I9").Value = 0.02980132450331128
$ws.Range("C10").Value = 0.1191860465116281
$ws.Range("D10").Value = 0
$ws.Range("G10").Value = 0.1154791154791162
$ws.Range("I10").Value = 0.01076158940397351
$ws.Range("J10").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.01017441860465117
$ws.Range("D12").Value = 0.007478632478632478
$ws.Range("G12").Value = 0
$ws.Range("I12").Value = 0.03559602649006623
$ws.Range("J12").Value = 0
$ws.Range("F13").Value = 0.05069124423963127
$ws.Range("H13").Value = 0
$ws.Range("B14").Value = 0.06689189189189185
$ws.Range("C14").Value = 0
$ws.Range("F14").Value = 0.1255760368663591
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0.3674157303370771
$ws.Range("C15").Value = 0.005813953488372093
$ws.Range("D15").Value = 0
$ws.Range("G15").Value = 0.01842751842751843
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0.0343992248062016
$ws.Range("D16").Value = 0.2948717948717942
$ws.Range("F16").Value = 0.01670506912442395
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0.04470198675496687
$ws.Range("J16").Value = 0
$ws.Range("F17").Value = 0.004032258064516129
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0.002247191011235955
$ws.Range("C18").Value = 0.0009689922480620155
$ws.Range("G18").Value = 0
$ws.Range("I18").Value = 0.00413907284768212
$ws.Range("J18").Value = 0
$ws.Range("B19").Value = 0.01283783783783784
$ws.Range("F19").Value = 0.09677419354838683
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0.003370786516853933
$ws.Range("F20").Value = 0.02361751152073733
$ws.Range("H20").Value = 0
$ws.Range("F21").Value = 0.04089861751152071
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0.01235955056179775
$ws.Range("G22").Value = 0.04750204750204763
$ws.Range("B23").Value = 0.0006756756756756757
$ws.Range("F23").Value = 0.126728110599078
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0.07752808988764047
$ws.Range("C24").Value = 0.003391472868217054
$ws.Range("G24").Value = 0
$ws.Range("I24").Value = 0.004966887417218543
$ws.Range("J24").Value = 0
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0.09108527131782961
$ws.Range("D32").Value = 0.3023504273504266
$ws.Range("G32").Value = 0.05159705159705176
$ws.Range("I32").Value = 0.06456953642384099
$ws.Range("J32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0.2349806201550391
$ws.Range("D33").Value = 0.009615384615384614
$ws.Range("G33").Value = 0.01556101556101556
$ws.Range("I33").Value = 0.06870860927152313
$ws.Range("J33").Value = 0
$ws.Range("F34").Value = 0.003456221198156682
$ws.Range("H34").Value = 0

# Remove now-unused "Joint regime area" rows 36-40 (full size / AF output fix)
$ws.Range("A36:K40").Delete()
